$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 21-32 with new chemical names and values
$ws.Range("A21").Value = "PFOA"
$ws.Range("B21").Value = 14.5
$ws.Range("C21").Value = 7

$ws.Range("B22").Value = 4
$ws.Range("C22").Value = 7.3

$ws.Range("B23").Value = 7.3
$ws.Range("C23").Value = 7.5

$ws.Range("B24").Value = 2.6

$ws.Range("B25").Value = 7
$ws.Range("C25").Value = 5.3

$ws.Range("B26").Value = 1.4
$ws.Range("C26").Value = 1.9

$ws.Range("B27").Value = 5.3
$ws.Range("C27").Value = 1.8

$ws.Range("B28").Value = 10.4
$ws.Range("C28").Value = 3.3

$ws.Range("B29").Value = 2
$ws.Range("C29").Value = 1.5

$ws.Range("B30").Value = 39
$ws.Range("C30").Value = 54.5

$ws.Range("B31").Value = 5.3
$ws.Range("C31").Value = 8.6

$ws.Range("A32").Value = "8:2 FTSA"
$ws.Range("B32").Value = 1.1
$ws.Range("C32").Value = 0.3

# Remove the now-obsolete last row (old row 33)
$ws.Rows.Item(33).Delete()
